# Adds new tracking/auth columns to several sheets (new attributes in the ERD),
# matching the commit "Added new attributes, hence updated ERD and Sample data."
#
# Order of the cell edits below matters: Excel appends brand-new shared strings
# to the shared-strings table in the order they are first entered, and the
# target workbook expects this exact order:
#   calories_date, oz_date, sleep_date, step_date, admin_password, user_password

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Calories sheet -> new trailing column C: "calories_date"
# ---------------------------------------------------------------------------
$wsCalories = $wb.Worksheets.Item("Calories")
$wsCalories.Cells.Item(1, 3).Value = "calories_date"

# ---------------------------------------------------------------------------
# 2) Hydration Level sheet -> new trailing column C: "oz_date"
# ---------------------------------------------------------------------------
$wsHydration = $wb.Worksheets.Item("Hydration Level")
$wsHydration.Cells.Item(1, 3).Value = "oz_date"

# ---------------------------------------------------------------------------
# 3) Sleep sheet -> new trailing column C: "sleep_date"
# ---------------------------------------------------------------------------
$wsSleep = $wb.Worksheets.Item("Sleep")
$wsSleep.Cells.Item(1, 3).Value = "sleep_date"

# ---------------------------------------------------------------------------
# 4) Steps sheet -> new trailing column D: "step_date"
# ---------------------------------------------------------------------------
$wsSteps = $wb.Worksheets.Item("Steps")
$wsSteps.Cells.Item(1, 4).Value = "step_date"
$wsSteps.Columns("D").ColumnWidth = 10.71

# ---------------------------------------------------------------------------
# 5) Admin info sheet -> new column C (inserted after email column B): "admin_password"
# ---------------------------------------------------------------------------
$wsAdmin = $wb.Worksheets.Item("Admin info")
$wsAdmin.Columns("C").Insert()
$wsAdmin.Cells.Item(1, 3).Value = "admin_password"
$wsAdmin.Columns("C").ColumnWidth = $wsAdmin.Columns("B").ColumnWidth

# ---------------------------------------------------------------------------
# 6) User sheet -> new column C (inserted after email column B): "user_password"
# ---------------------------------------------------------------------------
$wsUser = $wb.Worksheets.Item("User")
$wsUser.Columns("C").Insert()
$wsUser.Cells.Item(1, 3).Value = "user_password"
$wsUser.Columns("C").ColumnWidth = $wsUser.Columns("B").ColumnWidth

# ---------------------------------------------------------------------------
# View-state: restore/update each sheet's selection, then leave "User" as the
# active sheet/tab (previously "Calories" / sheet index 8 was active).
# ---------------------------------------------------------------------------
$wsHome = $wb.Worksheets.Item("Home")
$wsHome.Activate()
$wsHome.Range("E1").Select()

$wsAdmin.Activate()
$wsAdmin.Range("C1").Select()

$wsSteps.Activate()
$wsSteps.Range("D3").Select()

$wsSleep.Activate()
$wsSleep.Range("C1").Select()

$wsHydration.Activate()
$wsHydration.Range("E10").Select()

$wsCalories.Activate()
$wsCalories.Range("C1").Select()

$wsUser.Activate()
$wsUser.Range("J6").Select()
